$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1452202.8
$ws.Range("I17").Value = 471
$ws.Range("J17").Value = 1518190.5
$ws.Range("K17").Value = 1413
$ws.Range("L17").Value = 4554571.5
$ws.Range("M17").Value = -1245
$ws.Range("N17").Value = -4554907.5

$ws.Range("H18").Value = 700
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H62").Value = 3297.1177
$ws.Range("I62").Value = 2222.7778
$ws.Range("J62").Value = 4505.75
$ws.Range("K62").Value = 2222.7778
$ws.Range("L62").Value = 4505.75
$ws.Range("M62").Value = -1598.7778
$ws.Range("N62").Value = -5753.75

$ws.Range("H65").Value = 3297.1177
$ws.Range("I65").Value = 2222.7778
$ws.Range("J65").Value = 4505.75
$ws.Range("K65").Value = 11113.889
$ws.Range("L65").Value = 22528.75
$ws.Range("M65").Value = -7993.888999999999
$ws.Range("N65").Value = -28768.75

$ws.Range("H92").Value = 125000520
$ws.Range("I92").Value = 125000520
$ws.Range("K92").Value = 125000520
$ws.Range("M92").Value = -124999272

$ws.Range("H96").Value = 50001116
$ws.Range("I96").Value = 62501010
$ws.Range("J96").Value = 1550
$ws.Range("K96").Value = 187503030
$ws.Range("L96").Value = 4650
$ws.Range("M96").Value = -187501657
$ws.Range("N96").Value = -7396

$ws.Range("H103").Value = 1250225
$ws.Range("I103").Value = 1250225
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 3750675
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -3750089
$ws.Range("N103").ClearContents()

$ws.Range("H135").Value = 33346646
$ws.Range("I135").Value = 1648.75
$ws.Range("J135").Value = 71455220
$ws.Range("K135").Value = 14838.75
$ws.Range("L135").Value = 643096980
$ws.Range("M135").Value = -12303.75
$ws.Range("N135").Value = -643102050

$ws.Range("H137").Value = 45996.74
$ws.Range("I137").Value = 2805.2856
$ws.Range("K137").Value = 8415.856800000001
$ws.Range("M137").Value = -5865.856800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26388.023
$ws.Range("I32").Value = 26182.643
$ws.Range("K32").Value = 26182.643
$ws.Range("M32").Value = -25895.643

$ws.Range("H45").Value = 2888.4375
$ws.Range("I45").Value = 4900.2
$ws.Range("J45").Value = 1974
$ws.Range("K45").Value = 4900.2
$ws.Range("L45").Value = 1974
$ws.Range("M45").Value = -4523.2
$ws.Range("N45").Value = -2728

$ws.Range("H132").Value = 14068.571
$ws.Range("I132").Value = 1871.3334
$ws.Range("K132").Value = 5614.0002
$ws.Range("M132").Value = -3084.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3312.6843
$ws.Range("I20").Value = 3665
$ws.Range("K20").Value = 3665
$ws.Range("M20").Value = -3418

$ws.Range("H99").Value = 2144.4443
$ws.Range("I99").Value = 1760
$ws.Range("K99").Value = 1760
$ws.Range("M99").Value = -262

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 303.2
$ws.Range("J22").Value = 500.5
$ws.Range("L22").Value = 500.5
$ws.Range("N22").Value = -1200.5

$ws.Range("H50").Value = 14000
$ws.Range("J50").Value = 14000
$ws.Range("L50").Value = 14000
$ws.Range("N50").Value = -15250

$ws.Range("H51").Value = 14989.5
$ws.Range("J51").Value = 14989.5
$ws.Range("L51").Value = 14989.5
$ws.Range("N51").Value = -16461.5

$ws.Range("H59").Value = 21578.947
$ws.Range("J59").Value = 21578.947
$ws.Range("L59").Value = 21578.947
$ws.Range("N59").Value = -23868.947

$ws.Range("H60").Value = 11911.333
$ws.Range("J60").Value = 12317.883
$ws.Range("L60").Value = 12317.883
$ws.Range("N60").Value = -13339.883

$ws.Range("H61").Value = 14989.5
$ws.Range("J61").Value = 14989.5
$ws.Range("L61").Value = 14989.5
$ws.Range("N61").Value = -15685.5

$ws.Range("H74").Value = 33118.875
$ws.Range("J74").Value = 33118.875
$ws.Range("L74").Value = 33118.875
$ws.Range("N74").Value = -34866.875

$ws.Range("H77").Value = 33118.875
$ws.Range("J77").Value = 33118.875
$ws.Range("L77").Value = 99356.625
$ws.Range("N77").Value = -108092.625

$ws.Range("H86").Value = 5959872.5
$ws.Range("I86").Value = 2637.25
$ws.Range("J86").Value = 13902853
$ws.Range("K86").Value = 2637.25
$ws.Range("L86").Value = 13902853
$ws.Range("M86").Value = -1514.25
$ws.Range("N86").Value = -13905099

$ws.Range("H89").Value = 5959872.5
$ws.Range("I89").Value = 2637.25
$ws.Range("J89").Value = 13902853
$ws.Range("K89").Value = 13186.25
$ws.Range("L89").Value = 69514265
$ws.Range("M89").Value = -7570.25
$ws.Range("N89").Value = -69525497

$ws.Range("H94").Value = 2490.8572
$ws.Range("I94").Value = 1507.5
$ws.Range("K94").Value = 1507.5
$ws.Range("M94").Value = -1056.5

$ws.Range("H99").Value = 18522636
$ws.Range("I99").Value = 3414.9333
$ws.Range("J99").Value = 41671660
$ws.Range("K99").Value = 3414.9333
$ws.Range("L99").Value = 41671660
$ws.Range("M99").Value = -1916.9333
$ws.Range("N99").Value = -41674656

$ws.Range("H126").Value = 18522636
$ws.Range("I126").Value = 3414.9333
$ws.Range("J126").Value = 41671660
$ws.Range("K126").Value = 10244.7999
$ws.Range("L126").Value = 125014980
$ws.Range("M126").Value = -7774.7999
$ws.Range("N126").Value = -125019920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 7941.385
$ws.Range("J107").Value = 267.25
$ws.Range("L107").Value = 801.75
$ws.Range("N107").Value = -4641.75

$ws.Range("H122").Value = 751.1429000000001
$ws.Range("I122").Value = 412.5
$ws.Range("J122").Value = 886.6
$ws.Range("K122").Value = 3712.5
$ws.Range("L122").Value = 7979.400000000001
$ws.Range("M122").Value = -1262.5
$ws.Range("N122").Value = -12879.4

$ws.Range("H131").Value = 759.35
$ws.Range("I131").Value = 633
$ws.Range("J131").Value = 763.25775
$ws.Range("K131").Value = 1899
$ws.Range("L131").Value = 2289.77325
$ws.Range("M131").Value = 3141
$ws.Range("N131").Value = -12369.77325

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 58.533333
$ws.Range("I2").Value = 50.545456
$ws.Range("J2").Value = 80.5
$ws.Range("K2").Value = 50.545456
$ws.Range("L2").Value = 80.5
$ws.Range("M2").Value = 62.454544
$ws.Range("N2").Value = -306.5

$ws.Range("H132").Value = 47677.65
$ws.Range("I132").Value = 47412.477
$ws.Range("K132").Value = 142237.431
$ws.Range("M132").Value = -139707.431

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 941.05884
$ws.Range("I46").Value = 674.75
$ws.Range("J46").Value = 1177.7778
$ws.Range("K46").Value = 674.75
$ws.Range("L46").Value = 1177.7778
$ws.Range("M46").Value = -486.75
$ws.Range("N46").Value = -1553.7778

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1780
$ws.Range("I96").Value = 1966.6666
$ws.Range("J96").Value = 1500
$ws.Range("K96").Value = 1966.6666
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = -593.6666
$ws.Range("N96").Value = -4246

$ws.Range("H132").Value = 3999.5
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -19058
